$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content/format in the old used range
$ws.Range("A1:E7").Clear()

# Write new cell values
$ws.Range("A1").Value = "Title A1.1"
$ws.Range("B1").Value = "Title A1.2"
$ws.Range("C1").Value = "Title A1.3"
$ws.Range("A2").Value = "Static Row 1"
$ws.Range("A3").Value = "Static Row 2"
$ws.Range("A4").Value = "Static Row 3"
$ws.Range("A5").Value = "Static Row 4"
$ws.Range("A6").Value = "Static Row 5"
$ws.Range("A7").Value = "Static Row 6"
$ws.Range("A9").Value = "Titulo Static"
$ws.Range("B9").Value = "Content Status"
$ws.Range("C9").Value = "Dato Static"
$ws.Range("D9").Value = "Titulo"
$ws.Range("E9").Value = "Fecha de Publicacion"
$ws.Range("A10").Value = "Grupo de Datos 1"
$ws.Range("B10").Value = "standard"
$ws.Range("C10").Value = "Dato grupo 1"
$ws.Range("D10").Value = "U.S. Robotics presenta hallazgo"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "10-01-2029"
$ws.Range("E10").Style = "Normal"
$ws.Range("A11").Value = "Grupo de Datos 1"
$ws.Range("B11").Value = "standard"
$ws.Range("C11").Value = "Dato grupo 1"
$ws.Range("D11").Value = "Se presenta el nuevo teléfono móvil en evento"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "05-04-2030"
$ws.Range("E11").Style = "Normal"
$ws.Range("A12").Value = "Grupo de Datos 1"
$ws.Range("B12").Value = "standard"
$ws.Range("C12").Value = "Dato grupo 1"
$ws.Range("D12").Value = "Se mejora la conducción autónoma de vehículos"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "23-05-2022"
$ws.Range("E12").Style = "Normal"
$ws.Range("A13").Value = "Grupo de Datos 1"
$ws.Range("B13").Value = "standard"
$ws.Range("C13").Value = "Dato grupo 1"
$ws.Range("D13").Value = "Fuccia OS sacude al mundo"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "10-10-2028"
$ws.Range("E13").Style = "Normal"
$ws.Range("A14").Value = "Grupo de Datos 2"
$ws.Range("B14").Value = "standard"
$ws.Range("C14").Value = "Dato grupo 2"
$ws.Range("D14").Value = "Tenemos campeona del mundial de volleiball"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "09-09-2024"
$ws.Range("E14").Style = "Normal"
$ws.Range("A15").Value = "Grupo de Datos 2"
$ws.Range("B15").Value = "standard"
$ws.Range("C15").Value = "Dato grupo 2"
$ws.Range("D15").Value = "Equipo veterano da un gran espectaculo"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "01-12-2023"
$ws.Range("E15").Style = "Normal"
$ws.Range("A17").Value = "Title A2"
$ws.Range("A18").Value = "Subtitle A2.1"
$ws.Range("A19").Value = "Static Row 1"
$ws.Range("A21").Value = "Subtitle A2.2"
$ws.Range("A22").Value = "Static Row 2"
$ws.Range("A24").Value = "Subtitle A2.3"
$ws.Range("A25").Value = "Static Row 3"
$ws.Range("A27").Value = "Static Text: Lorem ipsum dolor sit amet, consectetur adipiscing elit. Quisque non laoreet mauris. Pellentesque habitant morbi tristique senectus et netus et malesuada fames ac turpis egestas. Curabitur vulputate bibendum nibh elementum pulvinar. Integer a leo in orci ultricies fermentum. Ut vitae velit et sapien congue accumsan sed tincidunt dui. Ut elementum imperdiet nunc, non hendrerit enim ultrices at. Sed rhoncus vehicula."
$ws.Range("A30").Value = "Title B1"
$ws.Range("A31").Value = "Title B1.1"
$ws.Range("B31").Value = "Title B1.2"
$ws.Range("C31").Value = "Title B1.3"
$ws.Range("D31").Value = "Title B1.4"
$ws.Range("E31").Value = "Title B1.5"
$ws.Range("F31").Value = "Title B1.6"
$ws.Range("G31").Value = "Title B1.7"
$ws.Range("A32").Value = "Subtitle 1"
$ws.Range("A33").Value = "Static Row 1"
$ws.Range("A34").Value = "Static Row 2"
$ws.Range("A35").Value = "Static Row 3"
$ws.Range("A36").Value = "Static Row 4"
$ws.Range("A37").Value = "Static Row 5"
$ws.Range("A38").Value = "Static Row 6"
$ws.Range("A39").Value = "Static Row 7"
$ws.Range("A40").Value = "Subtitle 2"
$ws.Range("A41").Value = "Static Row 1"

# Re-apply header style (bold, centered/top-aligned, thin border) to A9:E9
$hdr = $ws.Range("A9:E9")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2
